{"js": "const replacements = [\n  [\"84\u00d724=\", \"93\u00d785=\"],\n  [\"41\u00d775=\", \"69\u00d715=\"],\n  [\"59\u00d769=\", \"31\u00d725=\"],\n  [\"18\u00d747=\", \"38\u00d750=\"],\n  [\"12\u00d757=\", \"99\u00d753=\"],\n  [\"53\u00d734=\", \"29\u00d786=\"],\n  [\"17\u00d725=\", \"27\u00d786=\"],\n  [\"38\u00d755=\", \"89\u00d717=\"],\n  [\"12\u00d720=\", \"33\u00d754=\"],\n  [\"46\u00d730=\", \"18\u00d732=\"],\n  [\"45\u00d720=\", \"29\u00d717=\"],\n  [\"28\u00d733=\", \"54\u00d758=\"],\n  [\"28\u00d761=\", \"64\u00d796=\"],\n  [\"25\u00d781=\", \"87\u00d746=\"],\n  [\"45\u00d747=\", \"34\u00d755=\"],\n  [\"50\u00d776=\", \"42\u00d712=\"],\n  [\"27\u00d716=\", \"46\u00d729=\"],\n  [\"29\u00d794=\", \"20\u00d742=\"],\n  [\"70\u00d722=\", \"58\u00d783=\"],\n  [\"29\u00d754=\", \"96\u00d798=\"],\n  [\"27\u00d733=\", \"69\u00d728=\"],\n  [\"88\u00d725=\", \"98\u00d721=\"],\n  [\"29\u00d779=\", \"21\u00d719=\"],\n  [\"61\u00d793=\", \"98\u00d752=\"],\n  [\"39\u00d782=\", \"91\u00d735=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"84\u00d724=\"; New=\"93\u00d785=\"},\n    @{Old=\"41\u00d775=\"; New=\"69\u00d715=\"},\n    @{Old=\"59\u00d769=\"; New=\"31\u00d725=\"},\n    @{Old=\"18\u00d747=\"; New=\"38\u00d750=\"},\n    @{Old=\"12\u00d757=\"; New=\"99\u00d753=\"},\n    @{Old=\"53\u00d734=\"; New=\"29\u00d786=\"},\n    @{Old=\"17\u00d725=\"; New=\"27\u00d786=\"},\n    @{Old=\"38\u00d755=\"; New=\"89\u00d717=\"},\n    @{Old=\"12\u00d720=\"; New=\"33\u00d754=\"},\n    @{Old=\"46\u00d730=\"; New=\"18\u00d732=\"},\n    @{Old=\"45\u00d720=\"; New=\"29\u00d717=\"},\n    @{Old=\"28\u00d733=\"; New=\"54\u00d758=\"},\n    @{Old=\"28\u00d761=\"; New=\"64\u00d796=\"},\n    @{Old=\"25\u00d781=\"; New=\"87\u00d746=\"},\n    @{Old=\"45\u00d747=\"; New=\"34\u00d755=\"},\n    @{Old=\"50\u00d776=\"; New=\"42\u00d712=\"},\n    @{Old=\"27\u00d716=\"; New=\"46\u00d729=\"},\n    @{Old=\"29\u00d794=\"; New=\"20\u00d742=\"},\n    @{Old=\"70\u00d722=\"; New=\"58\u00d783=\"},\n    @{Old=\"29\u00d754=\"; New=\"96\u00d798=\"},\n    @{Old=\"27\u00d733=\"; New=\"69\u00d728=\"},\n    @{Old=\"88\u00d725=\"; New=\"98\u00d721=\"},\n    @{Old=\"29\u00d779=\"; New=\"21\u00d719=\"},\n    @{Old=\"61\u00d793=\"; New=\"98\u00d752=\"},\n    @{Old=\"39\u00d782=\"; New=\"91\u00d735=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute(\n        $r.Old,    # FindText\n        $false,    # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $r.New,    # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n}\n"}
